$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Booking")

# Update the Check-in / Check-out date test data used by the searchForStay test case.
# The leading apostrophe keeps the values as plain text (matching the existing
# text-as-date formatting of these cells) instead of being parsed into date serials.
$ws.Range("C2").Value = "'2023-10-15"
$ws.Range("D2").Value = "'2023-12-20"

# Move the active cell selection to C2, as recorded in the saved sheet view.
$ws.Range("C2").Select()
